$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of the original Excel date serial numbers (stored in column B, formatted
# as yyyy-mm-dd/custom date styles) to the literal text the cells should hold
# after the edit (dd/mm/yyyy, always the 1st of January).
$dates = @{
    40179 = "01/01/2010"
    40544 = "01/01/2011"
    40909 = "01/01/2012"
    41275 = "01/01/2013"
    41640 = "01/01/2014"
    42005 = "01/01/2015"
    42370 = "01/01/2016"
    42736 = "01/01/2017"
    43101 = "01/01/2018"
    43466 = "01/01/2019"
    43831 = "01/01/2020"
    44197 = "01/01/2021"
}

$lastRow = 37

for ($row = 2; $row -le $lastRow; $row++) {
    $bCell = $ws.Cells.Item($row, 2)
    $serial = [int]$bCell.Value2()
    $text = $dates[$serial]
    $bCell.Value = "'" + $text

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = $cCell.Value() * 100
}

# Reset the date column back to the default ("Normal") cell style now that the
# values are plain text, dropping the custom date number format entirely.
$ws.Range("B2:B" + $lastRow).Style = "Normal"
